$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "NA" values under duplicate_image_filename (column E) for rows 2-21
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
